# Apply the "want-to-go" counter bumps (展览 sheet) and the removal of the
# expired "上饶·万力时代次元企划嘉年华" listing plus counter bumps (全部类型 sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (exhibitions): only the F column (想去人数 / want-to-go count)
# changes for a handful of rows.
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Cells.Item(2, 6).Value = 55
$wsExpo.Cells.Item(4, 6).Value = 152
$wsExpo.Cells.Item(6, 6).Value = 5108
$wsExpo.Cells.Item(7, 6).Value = 110
$wsExpo.Cells.Item(8, 6).Value = 5285
$wsExpo.Cells.Item(10, 6).Value = 1341
$wsExpo.Cells.Item(11, 6).Value = 100

# ---------------------------------------------------------------------
# Sheet "全部类型" (all types): the expired free-entry listing in row 2
# ("上饶·万力时代次元企划嘉年华（免费展）", 2024-10-19) is removed, shifting
# every following row up by one. After the shift, re-number the sequence
# column (A) and apply the same counter bumps as above (rows shift by -1).
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Rows.Item(2).Delete()

$lastRow = $wsAll.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $wsAll.Cells.Item($r, 1).Value = $r - 1
}

$wsAll.Cells.Item(2, 6).Value = 55
$wsAll.Cells.Item(4, 6).Value = 152
$wsAll.Cells.Item(7, 6).Value = 5108
$wsAll.Cells.Item(8, 6).Value = 110
$wsAll.Cells.Item(9, 6).Value = 5285
$wsAll.Cells.Item(11, 6).Value = 1341
$wsAll.Cells.Item(12, 6).Value = 100
